$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a serial date value (46060 => 2026-02-07)
# that was bumped by one day (46061 => 2026-02-08) for every data row
# (rows 2 through 232).
$ws.Range("C2:C232").Value = 46061
